$d = $word.ActiveDocument

# Update the date heading (first paragraph)
$dateRange = $d.Paragraphs.Item(1).Range
$dateRange.Find.Execute("2025-09-21 Sunday", $false, $false, $false, $false, $false, $true, 1, $false, "2025-09-22 Monday", 2) | Out-Null

# Update each table cell with its new division problem
# (processed in an order that avoids transient text collisions)
$t = $d.Tables.Item(1)
$cell = $t.Cell(1, 1)
$cell.Range.Find.Execute("36÷8=4, 4", $false, $false, $false, $false, $false, $true, 1, $false, "13÷7=1, 6", 2) | Out-Null
$cell = $t.Cell(1, 2)
$cell.Range.Find.Execute("49÷3=16, 1", $false, $false, $false, $false, $false, $true, 1, $false, "96÷4=24, 0", 2) | Out-Null
$cell = $t.Cell(17, 5)
$cell.Range.Find.Execute("91÷7=13, 0", $false, $false, $false, $false, $false, $true, 1, $false, "97÷2=48, 1", 2) | Out-Null
$cell = $t.Cell(1, 3)
$cell.Range.Find.Execute("45÷4=11, 1", $false, $false, $false, $false, $false, $true, 1, $false, "91÷7=13, 0", 2) | Out-Null
$cell = $t.Cell(1, 4)
$cell.Range.Find.Execute("20÷9=2, 2", $false, $false, $false, $false, $false, $true, 1, $false, "55÷2=27, 1", 2) | Out-Null
$cell = $t.Cell(1, 5)
$cell.Range.Find.Execute("38÷2=19, 0", $false, $false, $false, $false, $false, $true, 1, $false, "83÷8=10, 3", 2) | Out-Null
$cell = $t.Cell(5, 1)
$cell.Range.Find.Execute("17÷4=4, 1", $false, $false, $false, $false, $false, $true, 1, $false, "23÷7=3, 2", 2) | Out-Null
$cell = $t.Cell(5, 2)
$cell.Range.Find.Execute("67÷4=16, 3", $false, $false, $false, $false, $false, $true, 1, $false, "66÷4=16, 2", 2) | Out-Null
$cell = $t.Cell(5, 3)
$cell.Range.Find.Execute("12÷7=1, 5", $false, $false, $false, $false, $false, $true, 1, $false, "30÷9=3, 3", 2) | Out-Null
$cell = $t.Cell(5, 4)
$cell.Range.Find.Execute("89÷9=9, 8", $false, $false, $false, $false, $false, $true, 1, $false, "86÷4=21, 2", 2) | Out-Null
$cell = $t.Cell(5, 5)
$cell.Range.Find.Execute("34÷2=17, 0", $false, $false, $false, $false, $false, $true, 1, $false, "24÷9=2, 6", 2) | Out-Null
$cell = $t.Cell(9, 1)
$cell.Range.Find.Execute("25÷2=12, 1", $false, $false, $false, $false, $false, $true, 1, $false, "76÷9=8, 4", 2) | Out-Null
$cell = $t.Cell(9, 2)
$cell.Range.Find.Execute("29÷9=3, 2", $false, $false, $false, $false, $false, $true, 1, $false, "84÷6=14, 0", 2) | Out-Null
$cell = $t.Cell(9, 3)
$cell.Range.Find.Execute("69÷7=9, 6", $false, $false, $false, $false, $false, $true, 1, $false, "79÷9=8, 7", 2) | Out-Null
$cell = $t.Cell(9, 4)
$cell.Range.Find.Execute("34÷7=4, 6", $false, $false, $false, $false, $false, $true, 1, $false, "92÷2=46, 0", 2) | Out-Null
$cell = $t.Cell(9, 5)
$cell.Range.Find.Execute("92÷6=15, 2", $false, $false, $false, $false, $false, $true, 1, $false, "83÷7=11, 6", 2) | Out-Null
$cell = $t.Cell(13, 1)
$cell.Range.Find.Execute("93÷9=10, 3", $false, $false, $false, $false, $false, $true, 1, $false, "81÷6=13, 3", 2) | Out-Null
$cell = $t.Cell(13, 2)
$cell.Range.Find.Execute("76÷2=38, 0", $false, $false, $false, $false, $false, $true, 1, $false, "90÷5=18, 0", 2) | Out-Null
$cell = $t.Cell(13, 3)
$cell.Range.Find.Execute("35÷7=5, 0", $false, $false, $false, $false, $false, $true, 1, $false, "66÷2=33, 0", 2) | Out-Null
$cell = $t.Cell(13, 4)
$cell.Range.Find.Execute("45÷3=15, 0", $false, $false, $false, $false, $false, $true, 1, $false, "38÷4=9, 2", 2) | Out-Null
$cell = $t.Cell(13, 5)
$cell.Range.Find.Execute("48÷6=8, 0", $false, $false, $false, $false, $false, $true, 1, $false, "11÷2=5, 1", 2) | Out-Null
$cell = $t.Cell(17, 1)
$cell.Range.Find.Execute("92÷8=11, 4", $false, $false, $false, $false, $false, $true, 1, $false, "55÷3=18, 1", 2) | Out-Null
$cell = $t.Cell(17, 2)
$cell.Range.Find.Execute("59÷2=29, 1", $false, $false, $false, $false, $false, $true, 1, $false, "94÷3=31, 1", 2) | Out-Null
$cell = $t.Cell(17, 3)
$cell.Range.Find.Execute("59÷9=6, 5", $false, $false, $false, $false, $false, $true, 1, $false, "20÷7=2, 6", 2) | Out-Null
$cell = $t.Cell(17, 4)
$cell.Range.Find.Execute("47÷5=9, 2", $false, $false, $false, $false, $false, $true, 1, $false, "74÷8=9, 2", 2) | Out-Null

Write-Output "Done"